$p = $ppt.ActivePresentation
$x = $p.ThisMethodDoesNotExist123("foo.thmx")
Write-Host "result: $x"
Write-Host "done"
